$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 227, pushing existing rows 227-239 down to 228-240
$ws.Rows.Item(227).Insert()

# Populate the newly inserted row 227 with the new weekly record
$ws.Cells.Item(227, 1).Value = 6
$ws.Cells.Item(227, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(227, 3).Value = "Metropolitana"
$ws.Cells.Item(227, 4).Value = 44826
$ws.Cells.Item(227, 5).Value = 13
$ws.Cells.Item(227, 6).Value = 100112029
$ws.Cells.Item(227, 7).Value = "Orégano"
$ws.Cells.Item(227, 8).Value = "Sin especificar"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 51
$ws.Cells.Item(227, 11).Value = 15000
$ws.Cells.Item(227, 12).Value = 16000
$ws.Cells.Item(227, 13).Value = 15451
$ws.Cells.Item(227, 14).Value = "$/docena de atados"
$ws.Cells.Item(227, 15).Value = "Región Metropolitana"
$ws.Cells.Item(227, 16).Value = 5150
$ws.Cells.Item(227, 17).Value = 3
$ws.Cells.Item(227, 18).Value = "Hortaliza"
